$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns keep their string representation (e.g. trailing zeros)
# for numeric-looking values, matching the original inline-string cell typing.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.329.25"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.55"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.73"
$ws.Range("E5").Value = "  +2.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.15"
$ws.Range("E6").Value = "  +4.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.89%  "

$ws.Range("E9").Value = "  +9.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.407"
$ws.Range("E10").Value = "  +5.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("E12").Value = "  +2.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.31"
$ws.Range("E13").Value = "  +6.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("E14").Value = "  +22.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.120.31"
$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.186.84"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.620.53"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("E18").Value = "  +5.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.90"
$ws.Range("E19").Value = "  +4.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "359.58"
$ws.Range("E20").Value = "  +4.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.37"
$ws.Range("E21").Value = "  +7.76%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.10"
$ws.Range("E23").Value = "  +3.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.68"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("E25").Value = "  +2.86%  "

$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.24"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  +3.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0966"
$ws.Range("E29").Value = "  +13.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  +9.81%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.41"
$ws.Range("E31").Value = "  -2.38%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.56"
$ws.Range("E34").Value = "  +5.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.39"
$ws.Range("E35").Value = "  +4.94%  "

$ws.Range("E36").Value = "  +4.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.69"
$ws.Range("E37").Value = "  +6.68%  "

$ws.Range("E38").Value = "  +5.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.83"
$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.69"
$ws.Range("E42").Value = "  +7.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.05"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.17"
$ws.Range("E44").Value = "  +4.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0620"
$ws.Range("E45").Value = "  +7.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.32"
$ws.Range("E46").Value = "  +2.21%  "

$ws.Range("E47").Value = "  +7.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.657"
$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("E49").Value = "  +6.94%  "

$ws.Range("E50").Value = "  +2.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.72"
$ws.Range("E51").Value = "  +3.62%  "
